$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.944.09'
$ws.Range("E2").Value = '  -1.97%  '

$ws.Range("D3").Value = '3.788.27'
$ws.Range("E3").Value = '  +2.57%  '

$ws.Range("E4").Value = '  +0.10%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '621.93'
$ws.Range("E5").Value = '  +3.78%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.41'
$ws.Range("E6").Value = '  -3.89%  '

$ws.Range("D7").Value = '3.787.97'
$ws.Range("E7").Value = '  +2.63%  '

$ws.Range("E8").Value = '  +0.05%  '

$ws.Range("E9").Value = '  +0.20%  '

$ws.Range("E10").Value = '  +4.19%  '

$ws.Range("E11").Value = '  -5.02%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '40.67'
$ws.Range("E13").Value = '  +1.95%  '

$ws.Range("E14").Value = '  +2.79%  '

$ws.Range("D15").Value = '4.424.31'
$ws.Range("E15").Value = '  +2.83%  '

$ws.Range("D16").Value = '3.789.64'
$ws.Range("E16").Value = '  +2.89%  '

$ws.Range("D17").Value = '69.972.15'
$ws.Range("E17").Value = '  -1.87%  '

$ws.Range("E18").Value = '  -0.14%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.60'
$ws.Range("E19").Value = '  +1.55%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.91'
$ws.Range("E20").Value = '  +0.05%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '509.33'
$ws.Range("E21").Value = '  -1.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.55'
$ws.Range("E22").Value = '  +3.66%  '

$ws.Range("E23").Value = '  -2.35%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.52'
$ws.Range("E24").Value = '  +3.82%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '87.42'
$ws.Range("E25").Value = '  +0.00%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.26'
$ws.Range("E26").Value = '  -1.24%  '

$ws.Range("B27").Value = 'RenderToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.13'
$ws.Range("E27").Value = '  +1.92%  '

$ws.Range("B28").Value = 'PEPE'
$ws.Range("C28").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000141'
$ws.Range("E28").Value = '  +27.91%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.03%  '

$ws.Range("E30").Value = '  -1.37%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.89'
$ws.Range("E31").Value = '  +4.30%  '

$ws.Range("E32").Value = '  -4.38%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.36'
$ws.Range("E33").Value = '  -1.37%  '

$ws.Range("E34").Value = '  -0.90%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.07'
$ws.Range("E36").Value = '  +5.27%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.21'
$ws.Range("E37").Value = '  +1.30%  '

$ws.Range("E38").Value = '  +4.51%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.333'
$ws.Range("E39").Value = '  -2.62%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.13'
$ws.Range("E40").Value = '  -1.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '50.95'
$ws.Range("E41").Value = '  -0.30%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '45.91'
$ws.Range("E42").Value = '  +1.84%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.73'
$ws.Range("E43").Value = '  -0.62%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '421.66'
$ws.Range("E44").Value = '  +3.01%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '3.038.66'
$ws.Range("E45").Value = '  -4.37%  '

$ws.Range("B46").Value = 'dogwifhat'
$ws.Range("C46").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.82'
$ws.Range("E46").Value = '  +1.51%  '

$ws.Range("E47").Value = '  -1.46%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.36'
$ws.Range("E48").Value = '  -3.13%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '138.27'
$ws.Range("E50").Value = '  +0.96%  '

$ws.Range("E51").Value = '  +1.36%  '
